$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$remarkCommon = "दिनांक __ रोजी रविवार असल्याने जमा झालेली रक्कम ही सोमवार दिनांक 03.12.2025 रोजी बँकेमध्ये भरणा करण्यात आली."
$remarkLast   = "दिनांक 02.12.2025 रोजी रविवार असल्याने जमा झालेली रक्कम ही सोमवार दिनांक 03.12.2025 रोजी बँकेमध्ये भरणा करण्यात आली."

$rows = @(
    @{ Row = 7;  A = "02-12-2025"; B = "010965012-Medha Sub Division Office Coll."; C = "Cash";   D = "2025-12-03"; E = 8770; F = $remarkCommon; G = "2025-12-23" },
    @{ Row = 8;  A = "02-12-2025"; B = "010965012-Medha Sub Division Office Coll."; C = "Cheque"; D = "";           E = 0;    F = $remarkCommon; G = "2025-12-23" },
    @{ Row = 9;  A = "02-12-2025"; B = "010965012-Medha Sub Division Office Coll."; C = "NEFT";   D = "";           E = 0;    F = $remarkCommon; G = "2025-12-23" },
    @{ Row = 10; A = "02-12-2025"; B = "010965012-Medha Sub Division Office Coll."; C = "Total";  D = "";           E = 0;    F = $remarkCommon; G = "2025-12-23" },
    @{ Row = 11; A = "02-12-2025"; B = "010965012-Medha Sub Division Office Coll."; C = "Cash";   D = "";           E = 0;    F = $remarkLast;   G = "2025-12-23" }
)

foreach ($r in $rows) {
    $cA = $ws.Cells.Item($r.Row, 1)
    $cA.NumberFormat = "@"
    $cA.Value = $r.A

    $cB = $ws.Cells.Item($r.Row, 2)
    $cB.NumberFormat = "@"
    $cB.Value = $r.B

    $cC = $ws.Cells.Item($r.Row, 3)
    $cC.NumberFormat = "@"
    $cC.Value = $r.C

    $cD = $ws.Cells.Item($r.Row, 4)
    $cD.NumberFormat = "@"
    $cD.Value = $r.D

    $ws.Cells.Item($r.Row, 5).Value = $r.E

    $cF = $ws.Cells.Item($r.Row, 6)
    $cF.NumberFormat = "@"
    $cF.Value = $r.F

    $cG = $ws.Cells.Item($r.Row, 7)
    $cG.NumberFormat = "@"
    $cG.Value = $r.G
}
